# Update countries & provincias Spain
#
# The "Illes Balears" data row (row 26) is removed from the Ciudades sheet;
# every row below it shifts up by one. This also leaves the "Illes Balears"
# shared string unused (it disappears from the workbook once nothing
# references it) while the already-present "Illes Balears*" entry is kept.
#
# Separately, the two neighbouring rows that used to read "Huesca" then
# "Huelva" are now ordered "Huelva" then "Huesca" (their underlying totals
# stay where they are - only the province labels swap).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire row for "Illes Balears" - everything beneath shifts up.
$ws.Range("A26").EntireRow.Delete()

# After the shift, the "Huesca"/"Huelva" pair now sits at rows 52-53;
# swap the two province labels so "Huelva" precedes "Huesca".
$ws.Range("A52").Value2 = "Huelva"
$ws.Range("A53").Value2 = "Huesca"
